$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.016718167572847
$ws.Range("C2").Value = 0.190278895237403

$ws.Range("B3").Value = 0.1003720825778019
$ws.Range("C3").Value = 0.1851143953343564

$ws.Range("B4").Value = 0.4392790356500113
$ws.Range("C4").Value = 0.143658757723793

$ws.Range("B5").Value = 0.9463933976029485
$ws.Range("C5").Value = 0.4239606987652852

$ws.Range("B6").Value = 0.7790492111245267
$ws.Range("C6").Value = 0.6048001395298276

$ws.Range("B7").Value = 0.4799275495545746
$ws.Range("C7").Value = 0.08973620266883328

$ws.Range("B8").Value = 0.006286856532096863
$ws.Range("C8").Value = 0.2229119110107422

$ws.Range("B9").Value = 0.08287183264522464
$ws.Range("C9").Value = 0.1406518752174287

$ws.Range("B10").Value = 0.8412007498430536
$ws.Range("C10").Value = 0.6575550125876612
